$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "43.053.70"
$ws.Range("D3").Value = "2.296.79"
$ws.Range("D5").Value = "'310.18"
$ws.Range("D6").Value = "'100.81"
$ws.Range("D7").Value = "'0.536"
$ws.Range("D9").Value = "'0.522"
$ws.Range("D10").Value = "'36.07"
$ws.Range("D11").Value = "'0.0823"
$ws.Range("D13").Value = "'7.11"
$ws.Range("D14").Value = "2.654.00"
$ws.Range("D16").Value = "2.295.51"
$ws.Range("D17").Value = "'0.806"
$ws.Range("D18").Value = "43.003.49"
$ws.Range("D20").Value = "0.0₃0921"
$ws.Range("D22").Value = "'68.46"
$ws.Range("D23").Value = "'240.19"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D27").Value = "'24.54"
$ws.Range("D28").Value = "'38.43"
$ws.Range("D29").Value = "'9.64"
$ws.Range("D30").Value = "'2.12"
$ws.Range("D31").Value = "'168.03"
$ws.Range("D32").Value = "'5.31"
$ws.Range("D35").Value = "'17.68"
$ws.Range("D36").Value = "'0.0738"
$ws.Range("D41").Value = "'4.20"
$ws.Range("D43").Value = "1.971.98"
$ws.Range("D45").Value = "'19.05"
$ws.Range("D48").Value = "'55.59"
$ws.Range("D50").Value = "2.525.20"

# Update Volume(1h) (column E) values, preserving the "  +x.xx%  " padding
$ws.Range("E2").Value = "  +2.89%  "
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("E6").Value = "  +6.22%  "
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +7.26%  "
$ws.Range("E10").Value = "  +4.11%  "
$ws.Range("E11").Value = "  +4.76%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  +7.56%  "
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("E15").Value = "  +4.99%  "
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").Value = "  +4.77%  "
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("E28").Value = "  +5.87%  "
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("E41").Value = "  +5.59%  "
$ws.Range("E42").Value = "  -3.22%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  +4.86%  "
$ws.Range("E49").Value = "  +16.34%  "
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("E51").Value = "  +2.27%  "

Write-Host "Cryptos list updated"